$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2307692307692308
$ws.Range("C2").Value = 0.4743589743589743
$ws.Range("J2").Value = 0.01282051282051282
$ws.Range("P2").Value = 0.1634615384615385
$ws.Range("S2").Value = 0.1185897435897436
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.06790123456790123
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7098765432098766
$ws.Range("S3").Value = 0.191358024691358
$ws.Range("J4").Value = 0.03773584905660377
$ws.Range("P4").Value = 0.5283018867924528
$ws.Range("S4").Value = 0.4339622641509434
$ws.Range("B6").Value = 0.03056768558951965
$ws.Range("D6").Value = 0.01746724890829694
$ws.Range("F6").Value = 0.07423580786026202
$ws.Range("J6").Value = 0.2925764192139738
$ws.Range("O6").Value = 0.02183406113537118
$ws.Range("Q6").Value = 0.1703056768558952
$ws.Range("R6").Value = 0.05240174672489083
$ws.Range("S6").Value = 0.3406113537117904
$ws.Range("B7").Value = 0.1045454545454545
$ws.Range("D7").Value = 0.01363636363636364
$ws.Range("F7").Value = 0.03636363636363636
$ws.Range("J7").Value = 0.15
$ws.Range("O7").Value = 0.02727272727272727
$ws.Range("Q7").Value = 0.1727272727272727
$ws.Range("R7").Value = 0.07727272727272727
$ws.Range("S7").Value = 0.4181818181818182
$ws.Range("B8").Value = 0.08888888888888889
$ws.Range("D8").Value = 0.0202020202020202
$ws.Range("E8").Value = 0.00202020202020202
$ws.Range("F8").Value = 0.05454545454545454
$ws.Range("J8").Value = 0.06666666666666667
$ws.Range("O8").Value = 0.01414141414141414
$ws.Range("Q8").Value = 0.1838383838383838
$ws.Range("R8").Value = 0.1111111111111111
$ws.Range("S8").Value = 0.4585858585858586
$ws.Range("B9").Value = 0.1185567010309278
$ws.Range("D9").Value = 0.02061855670103093
$ws.Range("F9").Value = 0.08247422680412371
$ws.Range("J9").Value = 0.05154639175257732
$ws.Range("O9").Value = 0.005154639175257732
$ws.Range("Q9").Value = 0.2216494845360825
$ws.Range("R9").Value = 0.06185567010309279
$ws.Range("S9").Value = 0.4381443298969072
$ws.Range("B10").Value = 0.1016702977487291
$ws.Range("D10").Value = 0.02396514161220044
$ws.Range("E10").Value = 0.0007262164124909223
$ws.Range("F10").Value = 0.05301379811183732
$ws.Range("J10").Value = 0.08496732026143791
$ws.Range("O10").Value = 0.0130718954248366
$ws.Range("Q10").Value = 0.2127814088598402
$ws.Range("R10").Value = 0.09658678286129267
$ws.Range("S10").Value = 0.4132171387073348
$ws.Range("G11").Value = 0.1803713527851459
$ws.Range("J11").Value = 0.07957559681697612
$ws.Range("K11").Value = 0.246684350132626
$ws.Range("L11").Value = 0.4748010610079575
$ws.Range("S11").Value = 0.01856763925729443
$ws.Range("G12").Value = 0.7098445595854922
$ws.Range("J12").Value = 0.1450777202072539
$ws.Range("K12").Value = 0.02590673575129534
$ws.Range("L12").Value = 0.07253886010362694
$ws.Range("S12").Value = 0.04663212435233161
$ws.Range("G13").Value = 0.6046511627906976
$ws.Range("J13").Value = 0.3488372093023256
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.04102564102564103
$ws.Range("H15").Value = 0.1230769230769231
$ws.Range("I15").Value = 0.07179487179487179
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.08205128205128205
$ws.Range("M15").Value = 0.01025641025641026
$ws.Range("O15").Value = 0.03076923076923077
$ws.Range("S15").Value = 0.241025641025641
$ws.Range("F16").Value = 0.02116402116402116
$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.06878306878306878
$ws.Range("J16").Value = 0.4391534391534391
$ws.Range("K16").Value = 0.1005291005291005
$ws.Range("M16").Value = 0.02645502645502645
$ws.Range("O16").Value = 0.0582010582010582
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.02994011976047904
$ws.Range("H17").Value = 0.1836327345309381
$ws.Range("I17").Value = 0.07984031936127745
$ws.Range("J17").Value = 0.4391217564870259
$ws.Range("K17").Value = 0.0718562874251497
$ws.Range("M17").Value = 0.02794411177644711
$ws.Range("N17").Value = 0.001996007984031936
$ws.Range("O17").Value = 0.04191616766467066
$ws.Range("S17").Value = 0.12375249500998
$ws.Range("F18").Value = 0.02202643171806168
$ws.Range("H18").Value = 0.1894273127753304
$ws.Range("I18").Value = 0.08370044052863436
$ws.Range("J18").Value = 0.4361233480176211
$ws.Range("K18").Value = 0.08370044052863436
$ws.Range("M18").Value = 0.02202643171806168
$ws.Range("O18").Value = 0.05286343612334802
$ws.Range("S18").Value = 0.1101321585903084
$ws.Range("F19").Value = 0.02108036890645586
$ws.Range("H19").Value = 0.2015810276679842
$ws.Range("I19").Value = 0.07312252964426877
$ws.Range("J19").Value = 0.3735177865612648
$ws.Range("K19").Value = 0.1231884057971015
$ws.Range("M19").Value = 0.01383399209486166
$ws.Range("O19").Value = 0.05072463768115942
$ws.Range("S19").Value = 0.1429512516469038
